$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D, shifting existing D:K data to F:M
$ws.Range("D:E").Insert()

# Copy number formats/styles from the (now shifted) F:M columns into the new D:E columns
$ws.Range("F7:M102").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D7").Value2 = 43464
$ws.Range("E7").Value2 = 43373
$ws.Range("D8").Value2 = 680000
$ws.Range("E8").Value2 = 660500
$ws.Range("D9").Value2 = 515500
$ws.Range("E9").Value2 = 499600
$ws.Range("D10").Value2 = 164500
$ws.Range("E10").Value2 = 160900
$ws.Range("D12").Value2 = "NA"
$ws.Range("E12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 0
$ws.Range("D14").Value2 = 2500
$ws.Range("E14").Value2 = 1100
$ws.Range("D15").Value2 = 0
$ws.Range("E15").Value2 = 0
$ws.Range("D17").Value2 = 630100
$ws.Range("E17").Value2 = 597100
$ws.Range("D18").Value2 = 49900
$ws.Range("E18").Value2 = 63400
$ws.Range("D20").Value2 = 100
$ws.Range("E20").Value2 = 1300
$ws.Range("D21").Value2 = 65200
$ws.Range("E21").Value2 = 78300
$ws.Range("D22").Value2 = 7100
$ws.Range("E22").Value2 = 6400
$ws.Range("D23").Value2 = 42900
$ws.Range("E23").Value2 = 58300
$ws.Range("D24").Value2 = 7800
$ws.Range("E24").Value2 = 10800
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 0
$ws.Range("D26").Value2 = 35100
$ws.Range("E26").Value2 = 47400
$ws.Range("D27").Value2 = 34900
$ws.Range("E27").Value2 = 47400
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 0
$ws.Range("D29").Value2 = 13500
$ws.Range("E29").Value2 = "NA"
$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("E31").Value2 = 0
$ws.Range("D32").Value2 = -100
$ws.Range("E32").Value2 = -1300
$ws.Range("D33").Value2 = 48400
$ws.Range("E33").Value2 = 47400
$ws.Range("D34").Value2 = 0
$ws.Range("E34").Value2 = 0
$ws.Range("D35").Value2 = 48400
$ws.Range("E35").Value2 = 47400
$ws.Range("D38").Value2 = 43464
$ws.Range("E38").Value2 = 43373
$ws.Range("D41").Value2 = 397200
$ws.Range("E41").Value2 = 545200
$ws.Range("D42").Value2 = 0
$ws.Range("E42").Value2 = 0
$ws.Range("D43").Value2 = 607400
$ws.Range("E43").Value2 = 519500
$ws.Range("D44").Value2 = 501300
$ws.Range("E44").Value2 = 396400
$ws.Range("D45").Value2 = 98900
$ws.Range("E45").Value2 = 81800
$ws.Range("D46").Value2 = 1604800
$ws.Range("E46").Value2 = 1543000
$ws.Range("D47").Value2 = 0
$ws.Range("E47").Value2 = 0
$ws.Range("D48").Value2 = 422200
$ws.Range("E48").Value2 = 398200
$ws.Range("D49").Value2 = 1129700
$ws.Range("E49").Value2 = 480700
$ws.Range("D50").Value2 = 0
$ws.Range("E50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("E51").Value2 = 0
$ws.Range("D52").Value2 = 61200
$ws.Range("E52").Value2 = 56800
$ws.Range("D53").Value2 = 0
$ws.Range("E53").Value2 = 0
$ws.Range("D54").Value2 = 3217900
$ws.Range("E54").Value2 = 2478600
$ws.Range("D57").Value2 = 281000
$ws.Range("E57").Value2 = 249900
$ws.Range("D58").Value2 = 3700
$ws.Range("E58").Value2 = 14100
$ws.Range("D59").Value2 = 266000
$ws.Range("E59").Value2 = 189800
$ws.Range("D60").Value2 = 550600
$ws.Range("E60").Value2 = 453800
$ws.Range("D61").Value2 = 1087300
$ws.Range("E61").Value2 = 599700
$ws.Range("D62").Value2 = 265700
$ws.Range("E62").Value2 = 222300
$ws.Range("D63").Value2 = 0
$ws.Range("E63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("E64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("E65").Value2 = 0
$ws.Range("D66").Value2 = 1908800
$ws.Range("E66").Value2 = 1280700
$ws.Range("D68").Value2 = 0
$ws.Range("E68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("E69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("E70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("E71").Value2 = 0
$ws.Range("D72").Value2 = 1439400
$ws.Range("E72").Value2 = 1398800
$ws.Range("D73").Value2 = 0
$ws.Range("E73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("E74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("E75").Value2 = 0
$ws.Range("D76").Value2 = 1309100
$ws.Range("E76").Value2 = 1197900
$ws.Range("D77").Value2 = 0
$ws.Range("E77").Value2 = 0
$ws.Range("D80").Value2 = 43464
$ws.Range("E80").Value2 = 43373
$ws.Range("D81").Value2 = 48400
$ws.Range("E81").Value2 = 47400
$ws.Range("D83").Value2 = 15200
$ws.Range("E83").Value2 = 13600
$ws.Range("D84").Value2 = 0
$ws.Range("E84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("E85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("E86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("E87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("E88").Value2 = 0
$ws.Range("D89").Value2 = 82400
$ws.Range("E89").Value2 = 58500
$ws.Range("D91").Value2 = -17200
$ws.Range("E91").Value2 = -20000
$ws.Range("D92").Value2 = 0
$ws.Range("E92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("E93").Value2 = 0
$ws.Range("D94").Value2 = -666800
$ws.Range("E94").Value2 = -19800
$ws.Range("D96").Value2 = -7500
$ws.Range("E96").Value2 = -7400
$ws.Range("D97").Value2 = 0
$ws.Range("E97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("E98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("E99").Value2 = 0
$ws.Range("D100").Value2 = 444400
$ws.Range("E100").Value2 = -1700
$ws.Range("D101").Value2 = -8000
$ws.Range("E101").Value2 = -4300
$ws.Range("D102").Value2 = -148000
$ws.Range("E102").Value2 = 32700

Write-Output "done"
